$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.613.51"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.798.93"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.45"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.14"
$ws.Range("E6").Value = "  -0.51%  "

# Row 7
$ws.Range("E7").Value = "  -0.65%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +6.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.12"
$ws.Range("E10").Value = "  -0.64%  "

# Row 11
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("E12").Value = "  -1.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.08"
$ws.Range("E13").Value = "  +3.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  +2.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.245.65"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.805.72"
$ws.Range("E16").Value = "  +0.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.940"
$ws.Range("E17").Value = "  -0.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.601.85"
$ws.Range("E18").Value = "  -0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.79"
$ws.Range("E19").Value = "  +3.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.18"
$ws.Range("E20").Value = "  +3.16%  "

# Row 21
$ws.Range("E21").Value = "  +1.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -0.38%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.36"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.14"
$ws.Range("E24").Value = "  -0.86%  "

# Row 25
$ws.Range("E25").Value = "  +0.74%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.11"
$ws.Range("E27").Value = "  -1.50%  "

# Row 28
$ws.Range("E28").Value = "  +1.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.36"
$ws.Range("E29").Value = "  +0.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.32"
$ws.Range("E30").Value = "  +7.47%  "

# Row 31
$ws.Range("E31").Value = "  +4.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.31"
$ws.Range("E32").Value = "  +9.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.12"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.68"
$ws.Range("E34").Value = "  +9.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0445"
$ws.Range("E35").Value = "  -5.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0853"
$ws.Range("E36").Value = "  +0.93%  "

# Row 37
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.84"
$ws.Range("E38").Value = "  -0.16%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("E39").Value = "  -1.81%  "

# Row 40
$ws.Range("E40").Value = "  -0.22%  "

# Row 41
$ws.Range("E41").Value = "  +0.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -4.57%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.99"
$ws.Range("E43").Value = "  +0.77%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.29"
$ws.Range("E44").Value = "  -0.48%  "

# Row 45
$ws.Range("E45").Value = "  -2.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.141.71"
$ws.Range("E46").Value = "  +2.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.44"
$ws.Range("E47").Value = "  +5.08%  "

# Row 48
$ws.Range("E48").Value = "  +6.90%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +11.00%  "

# Row 50
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.916"
$ws.Range("E50").Value = "  -4.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.219"
$ws.Range("E51").Value = "  +15.58%  "
